$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 5027.7144
$ws.Range("I19").Value = 1735.9048
$ws.Range("J19").Value = 9965.429
$ws.Range("K19").Value = 1735.9048
$ws.Range("L19").Value = 9965.429
$ws.Range("M19").Value = -1560.9048
$ws.Range("N19").Value = -10315.429

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1666.3334
$ws.Range("J38").Value = 5759
$ws.Range("L38").Value = 17277
$ws.Range("N38").Value = -18021

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 462.94116
$ws.Range("I53").Value = 557.3333
$ws.Range("J53").Value = 356.75
$ws.Range("K53").Value = 557.3333
$ws.Range("L53").Value = 356.75
$ws.Range("M53").Value = 79.66669999999999
$ws.Range("N53").Value = -1630.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4166.222
$ws.Range("I74").Value = 4166.222
$ws.Range("K74").Value = 4166.222
$ws.Range("M74").Value = -3230.222

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 4166.222
$ws.Range("I77").Value = 4166.222
$ws.Range("K77").Value = 20831.11
$ws.Range("M77").Value = -16151.11

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 5637.5654
$ws.Range("I86").Value = 7928.375
$ws.Range("J86").Value = 4415.8
$ws.Range("K86").Value = 7928.375
$ws.Range("L86").Value = 4415.8
$ws.Range("M86").Value = -6805.375
$ws.Range("N86").Value = -6661.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 5637.5654
$ws.Range("I89").Value = 7928.375
$ws.Range("J89").Value = 4415.8
$ws.Range("K89").Value = 39641.875
$ws.Range("L89").Value = 22079
$ws.Range("M89").Value = -34025.875
$ws.Range("N89").Value = -33311

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 45284
$ws.Range("J113").Value = 3732.4119
$ws.Range("L113").Value = 3732.4119
$ws.Range("N113").Value = -10240.4119

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 7365.143
$ws.Range("I137").Value = 7365.143
$ws.Range("K137").Value = 22095.429
$ws.Range("M137").Value = -19545.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 7842.5386
$ws.Range("I4").Value = 165.4
$ws.Range("K4").Value = 165.4
$ws.Range("M4").Value = -49.40000000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 387.5
$ws.Range("I5").Value = 200.16667
$ws.Range("J5").Value = 949.5
$ws.Range("K5").Value = 200.16667
$ws.Range("L5").Value = 949.5
$ws.Range("M5").Value = -88.16667000000001
$ws.Range("N5").Value = -1173.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6948.5415
$ws.Range("I45").Value = 7526
$ws.Range("J45").Value = 6266.091
$ws.Range("K45").Value = 7526
$ws.Range("L45").Value = 6266.091
$ws.Range("M45").Value = -7149
$ws.Range("N45").Value = -7020.091

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2087.739
$ws.Range("I74").Value = 1433.3846
$ws.Range("K74").Value = 1433.3846
$ws.Range("M74").Value = -559.3846000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2087.739
$ws.Range("I77").Value = 1433.3846
$ws.Range("K77").Value = 7166.923000000001
$ws.Range("M77").Value = -2798.923000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 119723.75
$ws.Range("J134").Value = 119723.75
$ws.Range("L134").Value = 119723.75
$ws.Range("N134").Value = -129863.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 56247.875
$ws.Range("J2").Value = 56247.875
$ws.Range("L2").Value = 56247.875
$ws.Range("N2").Value = -56473.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 387.5
$ws.Range("I4").Value = 200.16667
$ws.Range("J4").Value = 949.5
$ws.Range("K4").Value = 200.16667
$ws.Range("L4").Value = 949.5
$ws.Range("M4").Value = -85.16667000000001
$ws.Range("N4").Value = -1179.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 19174
$ws.Range("J6").Value = 19174
$ws.Range("L6").Value = 19174
$ws.Range("N6").Value = -19400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5551.643
$ws.Range("I20").Value = 6098.5
$ws.Range("J20").Value = 4184.5
$ws.Range("K20").Value = 6098.5
$ws.Range("L20").Value = 4184.5
$ws.Range("M20").Value = -5851.5
$ws.Range("N20").Value = -4678.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 9799
$ws.Range("J103").Value = 9799
$ws.Range("L103").Value = 9799
$ws.Range("N103").Value = -12143

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H116").Value = 131499.5
$ws.Range("J116").Value = 131499.5
$ws.Range("L116").Value = 131499.5
$ws.Range("N116").Value = -140677.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 49999.375
$ws.Range("J140").Value = 49999.375
$ws.Range("L140").Value = 49999.375
$ws.Range("N140").Value = -60359.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 209.57143
$ws.Range("I7").Value = 127.625
$ws.Range("K7").Value = 127.625
$ws.Range("M7").Value = -14.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 13250.444
$ws.Range("I22").Value = 16953.428
$ws.Range("J22").Value = 290
$ws.Range("K22").Value = 16953.428
$ws.Range("L22").Value = 290
$ws.Range("M22").Value = -16603.428
$ws.Range("N22").Value = -990

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3292.276
$ws.Range("I31").Value = 1705.4117
$ws.Range("J31").Value = 5540.3335
$ws.Range("K31").Value = 1705.4117
$ws.Range("L31").Value = 5540.3335
$ws.Range("M31").Value = -1410.4117
$ws.Range("N31").Value = -6130.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3292.276
$ws.Range("I34").Value = 1705.4117
$ws.Range("J34").Value = 5540.3335
$ws.Range("K34").Value = 1705.4117
$ws.Range("L34").Value = 5540.3335
$ws.Range("M34").Value = -1503.4117
$ws.Range("N34").Value = -5944.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1917.4667
$ws.Range("I58").Value = 1575.4
$ws.Range("J58").Value = 2601.6
$ws.Range("K58").Value = 1575.4
$ws.Range("L58").Value = 2601.6
$ws.Range("M58").Value = -1372.4
$ws.Range("N58").Value = -3007.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 849.4666999999999
$ws.Range("I107").Value = 749.75
$ws.Range("K107").Value = 749.75
$ws.Range("M107").Value = 1170.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H118").Value = 38382
$ws.Range("J118").Value = 38382
$ws.Range("L118").Value = 38382
$ws.Range("N118").Value = -41696

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H121").Value = 64075.332
$ws.Range("J121").Value = 64075.332
$ws.Range("L121").Value = 64075.332
$ws.Range("N121").Value = -66695.33199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H125").Value = 30194.25
$ws.Range("J125").Value = 30194.25
$ws.Range("L125").Value = 30194.25
$ws.Range("N125").Value = -35114.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1917.4667
$ws.Range("I136").Value = 1575.4
$ws.Range("J136").Value = 2601.6
$ws.Range("K136").Value = 4726.200000000001
$ws.Range("L136").Value = 7804.799999999999
$ws.Range("M136").Value = -2176.200000000001
$ws.Range("N136").Value = -12904.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 499.6154
$ws.Range("J5").Value = 90
$ws.Range("L5").Value = 270
$ws.Range("N5").Value = -494

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 89.5
$ws.Range("J12").Value = 81.666664
$ws.Range("L12").Value = 244.999992
$ws.Range("N12").Value = -590.999992

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 600.44446
$ws.Range("I14").Value = 600.44446
$ws.Range("K14").Value = 1801.33338
$ws.Range("M14").Value = -1628.33338

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1501.6
$ws.Range("I98").Value = 1564.5
$ws.Range("J98").Value = 1250
$ws.Range("K98").Value = 4693.5
$ws.Range("L98").Value = 3750
$ws.Range("M98").Value = -3195.5
$ws.Range("N98").Value = -6746

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 9999
$ws.Range("J120").Value = 9999
$ws.Range("L120").Value = 29997
$ws.Range("N120").Value = -39673

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1015324.75
$ws.Range("I131").Value = 1471267.2
$ws.Range("J131").Value = 103439.7
$ws.Range("K131").Value = 4413801.6
$ws.Range("L131").Value = 310319.1
$ws.Range("M131").Value = -4408761.6
$ws.Range("N131").Value = -320399.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 2711.1904
$ws.Range("I134").Value = 1683.4375
$ws.Range("J134").Value = 6000
$ws.Range("K134").Value = 5050.3125
$ws.Range("L134").Value = 18000
$ws.Range("M134").Value = 19.6875
$ws.Range("N134").Value = -28140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 499.6154
$ws.Range("J135").Value = 90
$ws.Range("L135").Value = 810
$ws.Range("N135").Value = -5880

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 1681.7693
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 4999.375
$ws.Range("J137").Value = 6166.6665
$ws.Range("L137").Value = 18499.9995
$ws.Range("N137").Value = -28699.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3238.85
$ws.Range("I102").Value = 3195.1836
$ws.Range("K102").Value = 3195.1836
$ws.Range("M102").Value = -1573.1836

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I22").Value = 996.5
$ws.Range("K22").Value = 996.5
$ws.Range("M22").Value = -701.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I27").Value = 996.5
$ws.Range("K27").Value = 996.5
$ws.Range("M27").Value = -889.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 40138.516
$ws.Range("I132").Value = 60856.277
$ws.Range("J132").Value = 11452.385
$ws.Range("K132").Value = 182568.831
$ws.Range("L132").Value = 34357.155
$ws.Range("M132").Value = -180038.831
$ws.Range("N132").Value = -39417.155

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3919888
$ws.Range("I136").Value = 7209586.5
$ws.Range("K136").Value = 21628759.5
$ws.Range("M136").Value = -21626209.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8269331
$ws.Range("I81").Value = 12990454
$ws.Range("K81").Value = 25980908
$ws.Range("M81").Value = -25979847

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 8269331
$ws.Range("I84").Value = 12990454
$ws.Range("K84").Value = 129904540
$ws.Range("M84").Value = -129899236

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3985.1428
$ws.Range("I132").Value = 4551.7144
$ws.Range("K132").Value = 13655.1432
$ws.Range("M132").Value = -11125.1432
